# "add prodi in mhs" — add a new "PRODI" column (I) to the mahasiswa_simple
# template sheet, to the right of the existing "NO BILLKEY" column (H).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell I1 = "PRODI" (goes into sharedStrings as a new unique entry).
# Leave formatting as the sheet's own default (column 9 is outside the
# explicitly-styled A:H block) so it picks up the same plain style already
# used for that column.
$ws.Range("I1").Value = "PRODI"

# Widen column I to fit the new header (~45 characters, matching column C's
# width class used elsewhere in the sheet).
$ws.Columns.Item(9).ColumnWidth = 44.285

# Reflect where the user ended up looking after adding the column.
$ws.Range("G19").Select()
